# Generate Report for Handoff
# The ff84103e-e205-4350-9f14-47d23e1a5a60.md record moved from
# "Handed back: in sync with en-US" to "Ready for handoff" with a new
# handoff timestamp and (for the per-language sheets) an error detail
# noting the handback file is stale.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/66ef41d5be6c5079429fab1b373a0296966b6109/e2e/ff84103e-e205-4350-9f14-47d23e1a5a60.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5c606d94b00445f33e3dc7dbf1e618a69f677c41/e2e/ff84103e-e205-4350-9f14-47d23e1a5a60.md."

# ---- Overview sheet: row 3 is the ff84103e record ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-09-07 13:25:24"

# ---- zh-cn sheet: row 3 is the ff84103e record ----
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("H3").Value = "2016-09-07 13:24:59"
$wsZhCn.Range("P3").Value = $errorDetail
$wsZhCn.Columns.Item(16).ColumnWidth = 39.17

# ---- de-de sheet: row 3 is the ff84103e record ----
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("H3").Value = "2016-09-07 13:25:24"
$wsDeDe.Range("P3").Value = $errorDetail
$wsDeDe.Columns.Item(16).ColumnWidth = 39.17
